# Regenerate merged AHB files
#
# 1. Rename the header row's "_old"/"_new" suffixed column headers to the
#    new release tags "_FV2410"/"_FV2504" (the "diff" column is untouched).
# 2. Turn the used range A1:U57 into an Excel Table ("Table1") that picks up
#    the (already renamed) header row as its column names.
# 3. Freeze the header row (freeze panes at A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the data range into an Excel Table, reusing the header row values
# as the table's column names.
$tableRange = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
